$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 229, pushing existing rows 229-265 down to 230-266.
$ws.Rows(229).Insert()

# Populate the newly inserted row 229 with its data (columns that stay
# constant across the whole sheet are also re-set here for completeness).
$ws.Range("A229").Value = 4
$ws.Range("B229").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C229").Value = "Los Lagos"
$ws.Range("D229").Value = 44951
$ws.Range("E229").Value = 10
$ws.Range("F229").Value = 100112009
$ws.Range("G229").Value = "Acelga"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 15
$ws.Range("K229").Value = 12000
$ws.Range("L229").Value = 12000
$ws.Range("M229").Value = 12000
$ws.Range("N229").Value = '$/docena de atados (12 kilos)'
$ws.Range("O229").Value = "Región de La Araucanía"
$ws.Range("P229").Value = 1000
$ws.Range("Q229").Value = 12
$ws.Range("R229").Value = "Hortaliza"
